$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Merge G1:H1 so the "Productos_Vendidos" header spans the new column too
$ws.Range("G1:H1").MergeCells = $true

# Merging recalculates the border styling on the merged range; restore the
# original bold/centered/bordered header style (same as A1) on G1:H1 so it
# matches the other header cells exactly.
$ws.Range("A1").Copy()
$ws.Range("G1:H1").PasteSpecial(-4122)
$excel.CutCopyMode = 0

# New "sum" sub-header in H2, matching the style of the existing sub-headers
$ws.Range("G2").Copy()
$ws.Range("H2").PasteSpecial(-4122)
$excel.CutCopyMode = 0
$ws.Range("H2").Value = "sum"

# New data values for the "Productos_Vendidos" sum column
$ws.Range("H4").Value = 520
$ws.Range("H5").Value = 413
